$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with fresh quote data.
# NumberFormat "@" keeps Excel from re-interpreting these literal strings
# as numbers/percentages; Style "Normal" afterwards drops the format
# override again so no cell style ends up attached (matches source, which
# has no "s" attribute on these cells).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "261.52"
Set-TextValue "D3" "26.45"
Set-TextValue "E3" "-3.41%"
Set-TextValue "D4" "4.717"
Set-TextValue "E4" "0.80%"
Set-TextValue "D5" "0.06158"
Set-TextValue "E5" "1.03%"
Set-TextValue "E6" "0.68%"
Set-TextValue "D7" "0.8512"
Set-TextValue "E7" "0.21%"
Set-TextValue "D8" "0.9124"
Set-TextValue "E8" "-1.23%"
Set-TextValue "D9" "0.1407"
Set-TextValue "E9" "0.43%"
Set-TextValue "D10" "0.05313"
Set-TextValue "E10" "8.71%"
Set-TextValue "D11" "0.07104"
Set-TextValue "E11" "0.13%"
Set-TextValue "D12" "0.03130"
Set-TextValue "E12" "1.73%"
Set-TextValue "D13" "0.09045"
Set-TextValue "E13" "-0.19%"
Set-TextValue "D14" "0.001540"
Set-TextValue "E14" "0.64%"
Set-TextValue "D15" "0.0006181"
Set-TextValue "E15" "1.26%"
Set-TextValue "D16" "0.005963"
Set-TextValue "E16" "-2.64%"
Set-TextValue "D17" "3.452"
Set-TextValue "E17" "0.03%"
Set-TextValue "E18" "0.85%"
Set-TextValue "E19" "1.07%"
Set-TextValue "D22" "4.086"
Set-TextValue "E22" "0.11%"
Set-TextValue "D23" "0.04225"
Set-TextValue "E23" "-0.16%"
Set-TextValue "D24" "0.001180"
Set-TextValue "E24" "-3.42%"
Set-TextValue "D25" "0.004048"
Set-TextValue "E25" "6.53%"
Set-TextValue "E26" "0.05%"
Set-TextValue "E27" "4.12%"
Set-TextValue "D40" "0.03989"
Set-TextValue "E40" "3.43%"
Set-TextValue "E41" "-0.04%"
Set-TextValue "D42" "0.004117"
Set-TextValue "E42" "1.11%"
Set-TextValue "D44" "0.01329"
Set-TextValue "E44" "-18.09%"
Set-TextValue "D45" "0.00005163"
Set-TextValue "E45" "0.22%"
Set-TextValue "E46" "0.05%"
Set-TextValue "D47" "0.02121"
Set-TextValue "D48" "0.2579"
Set-TextValue "E48" "90.34%"
Set-TextValue "E49" "0.05%"
Set-TextValue "E50" "0.05%"
